$wb = $excel.ActiveWorkbook

# "Tabelle2" (sheet1.xml) is the active/selected sheet in this workbook
$ws = $wb.Worksheets.Item("Tabelle2")
$ws.Select()

# Second check on the analyse_wh_word function: fill in newly reviewed
# Deviation_from_Biber / Precision notes for wh-questions and "that" complementation,
# plus a couple of "appears to work well" notes for wh-relatives features.
$ws.Range("D14").Value = "appears to work well"
$ws.Range("D34").Value = "appears to work well"
$ws.Range("D35").Value = "appears to work well"
$ws.Range("D23").Value = "serious problems here, catches unwanted stuff (`"there's nothing good that can come from it`", `"I'm sure that's a …`"), but also ignores some relevant examples without that (`"I am glad you liked it`")"
$ws.Range("C14").Value = "Biber excludes contracted auxiliaries here. I don't see why and our code currently does not exclude them"
$ws.Range("D22").Value = "I can't tell because the current taggeer will not tag the `"that`"s in my example sentences as WH…"

# Restore the reviewer's scroll position / selection on Tabelle2
$av = $excel.ActiveWindow
$av.ScrollRow = 19
$av.ScrollColumn = 2
$ws.Range("D22").Select()

$wb.Save()
